# Update cryptos list: refreshed prices and 1h volume percentage changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new Price values below read as plain numbers; reassert each cell as
# Text first so Excel keeps the exact string instead of converting it to a
# floating-point number (matching how these columns were already stored).
foreach ($addr in @("D5","D8","D10","D11","D14","D16","D17","D18","D19","D21","D23","D24","D25","D27","D32","D33","D34","D37","D40","D41","D43","D44","D46","D50","D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '34.082.10'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '1.778.50'
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("D5").Value = '225.13'
$ws.Range("E5").Value = '  -0.55%  '
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").Value = '31.67'
$ws.Range("E8").Value = '  -1.15%  '
$ws.Range("E9").Value = '  -0.43%  '
$ws.Range("D10").Value = '0.0684'
$ws.Range("E10").Value = '  +0.56%  '
$ws.Range("D11").Value = '0.0947'
$ws.Range("E11").Value = '  +1.00%  '
$ws.Range("D12").Value = '2.035.17'
$ws.Range("E12").Value = '  -0.21%  '
$ws.Range("D13").Value = '1.779.61'
$ws.Range("E13").Value = '  +0.11%  '
$ws.Range("D14").Value = '10.86'
$ws.Range("E14").Value = '  -2.56%  '
$ws.Range("D15").Value = '34.125.21'
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("D16").Value = '0.620'
$ws.Range("E16").Value = '  +0.93%  '
$ws.Range("D17").Value = '4.18'
$ws.Range("E17").Value = '  +0.28%  '
$ws.Range("D18").Value = '67.51'
$ws.Range("E18").Value = '  +0.01%  '
$ws.Range("D19").Value = '243.97'
$ws.Range("E19").Value = '  +0.85%  '
$ws.Range("D20").Value = '0.0₃0785'
$ws.Range("E20").Value = '  +1.91%  '
$ws.Range("D21").Value = '10.99'
$ws.Range("E21").Value = '  +3.47%  '
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("D23").Value = '4.09'
$ws.Range("E23").Value = '  +0.42%  '
$ws.Range("D24").Value = '2.04'
$ws.Range("E24").Value = '  -1.23%  '
$ws.Range("D25").Value = '160.80'
$ws.Range("E25").Value = '  -0.61%  '
$ws.Range("E26").Value = '  -0.31%  '
$ws.Range("D27").Value = '16.19'
$ws.Range("E27").Value = '  +0.25%  '
$ws.Range("E28").Value = '  +0.82%  '
$ws.Range("E29").Value = '  +0.21%  '
$ws.Range("E30").Value = '  -0.46%  '
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("D32").Value = '3.71'
$ws.Range("E32").Value = '  +2.56%  '
$ws.Range("D33").Value = '3.69'
$ws.Range("E33").Value = '  +3.79%  '
$ws.Range("D34").Value = '1.79'
$ws.Range("E34").Value = '  -2.24%  '
$ws.Range("D35").Value = '1.439.57'
$ws.Range("E35").Value = '  +3.31%  '
$ws.Range("E36").Value = '  +1.30%  '
$ws.Range("D37").Value = '2.44'
$ws.Range("E37").Value = '  +4.75%  '
$ws.Range("E38").Value = '  +1.16%  '
$ws.Range("E39").Value = '  +0.55%  '
$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").Value = '80.20'
$ws.Range("E40").Value = '  +0.81%  '
$ws.Range("B41").Value = 'HuobiToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D41").Value = '2.37'
$ws.Range("E41").Value = '  +0.60%  '
$ws.Range("E42").Value = '  +1.83%  '
$ws.Range("D43").Value = '0.913'
$ws.Range("E43").Value = '  -0.52%  '
$ws.Range("D44").Value = '13.57'
$ws.Range("E44").Value = '  -0.75%  '
$ws.Range("E45").Value = '  +1.13%  '
$ws.Range("D46").Value = '6.04'
$ws.Range("E46").Value = '  +2.37%  '
$ws.Range("E47").Value = '  -0.62%  '
$ws.Range("D48").Value = '1.936.99'
$ws.Range("E48").Value = '  -0.23%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0131'
$ws.Range("E49").Value = '  -7.10%  '
$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  +0.19%  '
$ws.Range("D51").Value = '103.77'
$ws.Range("E51").Value = '  -3.11%  '
